# Weekly update: add a new week of data (rows 338-340) for Repollo prices,
# shifting the existing historical rows down, plus a few small date/origin
# corrections on rows that were re-dated as part of this weekly refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows before the current row 338, pushing the existing
# rows 338:355 down to 341:358.
$ws.Range("A338:R340").EntireRow.Insert()

# --- Populate the 3 newly inserted rows (new week: 2021-11-16 / 44516) ---

# Row 338
$ws.Cells.Item(338, 1).Value()  = 9
$ws.Cells.Item(338, 2).Value()  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(338, 3).Value()  = "Metropolitana"
$ws.Cells.Item(338, 4).Value()  = 44516
$ws.Cells.Item(338, 5).Value()  = 13
$ws.Cells.Item(338, 6).Value()  = 100112006
$ws.Cells.Item(338, 7).Value()  = "Repollo"
$ws.Cells.Item(338, 8).Value()  = "Crespo record"
$ws.Cells.Item(338, 9).Value()  = "Primera"
$ws.Cells.Item(338, 10).Value() = 5200
$ws.Cells.Item(338, 11).Value() = 600
$ws.Cells.Item(338, 12).Value() = 700
$ws.Cells.Item(338, 13).Value() = 650
$ws.Cells.Item(338, 14).Value() = "`$/unidad"
$ws.Cells.Item(338, 15).Value() = "Región Metropolitana"
$ws.Cells.Item(338, 16).Value() = 650
$ws.Cells.Item(338, 17).Value() = 1
$ws.Cells.Item(338, 18).Value() = "Hortaliza"

# Row 339
$ws.Cells.Item(339, 1).Value()  = 9
$ws.Cells.Item(339, 2).Value()  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(339, 3).Value()  = "Metropolitana"
$ws.Cells.Item(339, 4).Value()  = 44516
$ws.Cells.Item(339, 5).Value()  = 13
$ws.Cells.Item(339, 6).Value()  = 100112006
$ws.Cells.Item(339, 7).Value()  = "Repollo"
$ws.Cells.Item(339, 8).Value()  = "Crespo record"
$ws.Cells.Item(339, 9).Value()  = "Segunda"
$ws.Cells.Item(339, 10).Value() = 1600
$ws.Cells.Item(339, 11).Value() = 500
$ws.Cells.Item(339, 12).Value() = 500
$ws.Cells.Item(339, 13).Value() = 500
$ws.Cells.Item(339, 14).Value() = "`$/unidad"
$ws.Cells.Item(339, 15).Value() = "Región Metropolitana"
$ws.Cells.Item(339, 16).Value() = 500
$ws.Cells.Item(339, 17).Value() = 1
$ws.Cells.Item(339, 18).Value() = "Hortaliza"

# Row 340
$ws.Cells.Item(340, 1).Value()  = 9
$ws.Cells.Item(340, 2).Value()  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(340, 3).Value()  = "Metropolitana"
$ws.Cells.Item(340, 4).Value()  = 44516
$ws.Cells.Item(340, 5).Value()  = 13
$ws.Cells.Item(340, 6).Value()  = 100112006
$ws.Cells.Item(340, 7).Value()  = "Repollo"
$ws.Cells.Item(340, 8).Value()  = "Morada(o)"
$ws.Cells.Item(340, 9).Value()  = "Primera"
$ws.Cells.Item(340, 10).Value() = 1600
$ws.Cells.Item(340, 11).Value() = 800
$ws.Cells.Item(340, 12).Value() = 900
$ws.Cells.Item(340, 13).Value() = 850
$ws.Cells.Item(340, 14).Value() = "`$/unidad"
$ws.Cells.Item(340, 15).Value() = "Región Metropolitana"
$ws.Cells.Item(340, 16).Value() = 850
$ws.Cells.Item(340, 17).Value() = 1
$ws.Cells.Item(340, 18).Value() = "Hortaliza"

# --- Small corrections on the shifted-down historical rows ---
# (dates / origin that were adjusted as part of this weekly refresh)

$ws.Cells.Item(344, 4).Value()  = 44509
$ws.Cells.Item(345, 4).Value()  = 44509
$ws.Cells.Item(349, 4).Value()  = 44421
$ws.Cells.Item(350, 4).Value()  = 44421
$ws.Cells.Item(353, 15).Value() = "Región de O'Higgins"
